$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Fix header text typos in row 1 ---
$ws.Range("C1").Value = "distance by car (km)"
$ws.Range("D1").Value = "time by car (hours)"

# --- 2. Give existing data cells (A2:E9) an explicit white fill (fillId=2), as in target ---
$ws.Range("A2:E9").Interior.ColorIndex = 2

# --- 3. Add new column F: "transatlantic flight equivalents (flights)" ---
# Header cell F1: same style as the rest of the header row (copy from E1)
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "transatlantic flight equivalents (flights)"

# Data cell F2: copy formatting (border/font/fill) from E2, then make it General number format
$ws.Range("E2").Copy()
$ws.Range("F2").PasteSpecial(-4122)
$ws.Range("F2").NumberFormat = "General"
$ws.Range("F2").Value = 0.02941

# Data cells F3:F9: copy formatting from E3, then make General number format
$ws.Range("E3").Copy()
$ws.Range("F3:F9").PasteSpecial(-4122)
$ws.Range("F3:F9").NumberFormat = "General"

$ws.Range("F3").Value = 0.066206
$ws.Range("F4").Value = 0.122151
$ws.Range("F5").Value = 0.092459
$ws.Range("F6").Value = 0.13175
$ws.Range("F7").Value = 0.080014
$ws.Range("F8").Value = 0.037342
$ws.Range("F9").Value = 0.036366

$excel.CutCopyMode = 0

# --- 4. Remove the trailing blank row 10 ---
$ws.Rows.Item(10).Delete()

# --- 5. Column widths: C and D change slightly, E/F share the same width ---
$ws.Columns.Item(3).ColumnWidth = 7.333333333333333
$ws.Columns.Item(4).ColumnWidth = 36.166666666666664
$ws.Columns.Item(6).ColumnWidth = 8.0

Write-Host "done"
